
# Workbook / sheet handles -------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rights_and_functions")

# A new script/trigger-template row is inserted above the existing row 57
# ("46_cre_view_fe_dataproc_last.sql"), pushing the remaining rows of the
# table down by one. This mirrors the xlsx diff which adds a brand-new row
# for "45_cre_table_frontend_in_trig.sql" / "template_cre_trigger_set_id.sql"
# right before that entry.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new trigger-generation entry.
$ws.Cells.Item(57, 2).Value = "45_cre_table_frontend_in_trig.sql"   # B57 SCRIPTNAME
$ws.Cells.Item(57, 3).Value = "template_cre_trigger_set_id.sql"     # C57 TEMPLATE
$ws.Cells.Item(57, 4).Value = "db2frontend_user"                    # D57 OWNER_USER
$ws.Cells.Item(57, 5).Value = "db2frontend_in"                      # E57 OWNER_SCHEMA
$ws.Cells.Item(57, 8).Value = "_fe"                                 # H57 TABLE_POSTFIX

# Move the active selection to where the author left it after the edit.
$ws.Activate()
$ws.Range("H58").Select()
